$wb = $excel.ActiveWorkbook

# --- Sheet 1: "概况" (Overview) -------------------------------------------
# Numeric stat cells: caliper reverted from 0.01 to 0.02, matched sample
# grew from 2648 to 2846 obs (1324/1324 -> 1423/1423 per group).
$ws1 = $wb.Worksheets.Item("概况")

$ws1.Range("B3").Value = 2846
$ws1.Range("B4").Value = 1423
$ws1.Range("B5").Value = 1423
$ws1.Range("B8").Value = 0.02

# --- Sheet 2: "平衡性汇总" (Balance summary) --------------------------------
# C/D columns hold post-match bias % and bias-reduction % as plain TEXT
# (not numbers) in the original workbook, so force text formatting before
# writing the new figures and restore the default style afterward so no
# visible formatting changes stick to the cells.
$ws2 = $wb.Worksheets.Item("平衡性汇总")

$textCells = @(
    @{ Addr = "C2"; Val = "-2.61" },
    @{ Addr = "D2"; Val = "106.4" },
    @{ Addr = "C3"; Val = "1.29" },
    @{ Addr = "D3"; Val = "96.7" },
    @{ Addr = "C4"; Val = "-3.86" },
    @{ Addr = "D4"; Val = "117.6" },
    @{ Addr = "C5"; Val = "7.23" },
    @{ Addr = "D5"; Val = "76.1" },
    @{ Addr = "C6"; Val = "-8.77" },
    @{ Addr = "D6"; Val = "-60.4" }
)

foreach ($cell in $textCells) {
    $rng = $ws2.Range($cell.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cell.Val
    $rng.Style = "Normal"
}
